# The sheet originally listed 13 categories in A1:A13, with two duplicate
# entries accidentally included ("Cartão de Crédito" in row 2, a near-dupe
# of "Cartão de crédico" that stayed in row 12, and "Outros" repeated in
# row 11). Remove the three redundant rows so the list shrinks from
# A1:A13 down to A1:A10, shifting the remaining rows up.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Cartão de Crédito" (duplicate credit-card category)
$ws.Rows(2).Delete()

# After the shift above, the old rows 11 and 12 ("Outros" duplicate and
# "Cartão de crédico") are now rows 10 and 11. Deleting row 10 twice
# removes both, leaving the old row 13 ("Investimentos") as the new row 10.
$ws.Rows(10).Delete()
$ws.Rows(10).Delete()
